# feat(dataset): add csv/xlsx thresholds
#
# Renames the default worksheet to "dataset" and appends six new
# "threshold_*" header columns (F1:K1) after the existing
# id/question/answer/contexts/ground_truth headers (A1:E1), matching the
# bare/plain formatting used across the whole header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: Sheet1 -> dataset
$ws.Name = "dataset"

# Full set of header values for row 1 (A1:K1)
$headers = @(
    "id",
    "question",
    "answer",
    "contexts",
    "ground_truth",
    "threshold_faithfulness",
    "threshold_answer_relevancy",
    "threshold_context_precision",
    "threshold_context_recall",
    "threshold_factual_correctness",
    "threshold_semantic_similarity"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}

# The original template bolded/boxed A1:E1; the new template keeps every
# header cell in the plain default style, so strip that formatting.
$ws.Range("A1:K1").ClearFormats()

# Reset page margins back to the regular Excel defaults (inches):
# left/right 0.75", top/bottom 1", header/footer 0.5"
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
